# Adds a new forecast date (2020-05-07) as a new column "AB" and a new
# observed date (2020-05-13) as a new row (40) to both the "cases" and
# "deaths" sheets, fills in the diagonal forecast values for the new
# column, fills in the previously-missing "Observed" value for
# 2020-04-29 (row 26, column B), and sets the new row's date label.

$wb = $excel.ActiveWorkbook

# New forecast values for column AB (rows 27-40), keyed by row number.
$caseValues = @{
    27 = 84920
    28 = 92411
    29 = 100756
    30 = 108573
    31 = 118002
    32 = 127991
    33 = 138144
    34 = 147911
    35 = 157099
    36 = 165135
    37 = 172733
    38 = 180048
    39 = 186892
    40 = 193040
}

$deathValues = @{
    27 = 5986
    28 = 6559
    29 = 7215
    30 = 7829
    31 = 8592
    32 = 9393
    33 = 10187
    34 = 10955
    35 = 11675
    36 = 12299
    37 = 12903
    38 = 13483
    39 = 14011
    40 = 14503
}

# The previously-missing "Observed" value for 2020-04-29 (row 26).
$observedB26 = @{
    "cases" = 26158
    "deaths" = 2247
}

$sheetValues = @{
    "cases" = $caseValues
    "deaths" = $deathValues
}

foreach ($sheetName in @("cases", "deaths")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New header date for the new forecast column AB (col 28). Force text
    # storage so Excel doesn't coerce the date-like string into a date
    # serial number (matches the existing header cells, which are text).
    $ws.Cells.Item(1, 28).NumberFormat = "@"
    $ws.Cells.Item(1, 28).Value = "2020-05-07"

    # Fill the new forecast values down column AB.
    $values = $sheetValues[$sheetName]
    foreach ($r in $values.Keys) {
        $ws.Cells.Item($r, 28).Value = $values[$r]
    }

    # Fill in the previously-missing "Observed" figure for 2020-04-29.
    $ws.Cells.Item(26, 2).Value = $observedB26[$sheetName]

    # New row 40 for the new observed date 2020-05-13 (column A), again
    # forcing text storage to match the rest of column A.
    $ws.Cells.Item(40, 1).NumberFormat = "@"
    $ws.Cells.Item(40, 1).Value = "2020-05-13"
}
